$wb = $excel.ActiveWorkbook

$timestamp = "2025-11-29 03:04:37"

# --- Главные ---
$ws = $wb.Worksheets.Item("Главные")

# Row 2
$ws.Range("C2").Value = 31
$ws.Range("D2").Value = 655
$ws.Range("E2").Value = 275
$ws.Range("F2").Value = 380
$ws.Range("G2").Value = 21.13
$ws.Range("H2").Value = 8.869999999999999
$ws.Range("I2").Value = 12.26
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 145

# Row 12
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 329
$ws.Range("E12").Value = 141
$ws.Range("F12").Value = 188
$ws.Range("G12").Value = 17.32
$ws.Range("H12").Value = 7.42
$ws.Range("I12").Value = 9.890000000000001
$ws.Range("J12").Value = 58
$ws.Range("K12").Value = 69

# Row 14
$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 288
$ws.Range("E14").Value = 151
$ws.Range("F14").Value = 137
$ws.Range("G14").Value = 13.71
$ws.Range("H14").Value = 7.19
$ws.Range("I14").Value = 6.52
$ws.Range("J14").Value = 68
$ws.Range("K14").Value = 56
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 5
$ws.Range("V14").Value = 12

# Row 16
$ws.Range("C16").Value = 29
$ws.Range("D16").Value = 531
$ws.Range("E16").Value = 262
$ws.Range("F16").Value = 269
$ws.Range("G16").Value = 18.31
$ws.Range("H16").Value = 9.029999999999999
$ws.Range("I16").Value = 9.279999999999999
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = 102
$ws.Range("W16").Value = 8

# Row 19
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 408
$ws.Range("E19").Value = 202
$ws.Range("F19").Value = 206
$ws.Range("G19").Value = 17.74
$ws.Range("H19").Value = 8.779999999999999
$ws.Range("I19").Value = 8.960000000000001
$ws.Range("J19").Value = 96
$ws.Range("K19").Value = 88

# Row 22
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 444
$ws.Range("E22").Value = 195
$ws.Range("F22").Value = 249
$ws.Range("G22").Value = 21.14
$ws.Range("H22").Value = 9.289999999999999
$ws.Range("I22").Value = 11.86
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 87

# Row 24
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 495
$ws.Range("E24").Value = 226
$ws.Range("F24").Value = 269
$ws.Range("G24").Value = 16.5
$ws.Range("H24").Value = 7.53
$ws.Range("J24").Value = 108
$ws.Range("K24").Value = 117
$ws.Range("L24").Value = 2
$ws.Range("M24").Value = 3
$ws.Range("V24").Value = 18

foreach ($r in 2..26) {
    $ws.Range("AA$r").Value = $timestamp
}

# --- Линейные ---
$ws = $wb.Worksheets.Item("Линейные")

# Row 2
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 360
$ws.Range("E2").Value = 158
$ws.Range("F2").Value = 202
$ws.Range("G2").Value = 18.95
$ws.Range("H2").Value = 8.32
$ws.Range("I2").Value = 10.63
$ws.Range("J2").Value = 69
$ws.Range("K2").Value = 76

# Row 8
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = 421
$ws.Range("E8").Value = 158
$ws.Range("F8").Value = 263
$ws.Range("G8").Value = 15.59
$ws.Range("H8").Value = 5.85
$ws.Range("I8").Value = 9.74
$ws.Range("J8").Value = 74
$ws.Range("K8").Value = 104

# Row 12
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = 422
$ws.Range("E12").Value = 202
$ws.Range("F12").Value = 220
$ws.Range("G12").Value = 17.58
$ws.Range("H12").Value = 8.42
$ws.Range("J12").Value = 91
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 4
$ws.Range("V12").Value = 16

# Row 15
$ws.Range("C15").Value = 24
$ws.Range("D15").Value = 455
$ws.Range("E15").Value = 235
$ws.Range("F15").Value = 220
$ws.Range("G15").Value = 18.96
$ws.Range("H15").Value = 9.789999999999999
$ws.Range("I15").Value = 9.17
$ws.Range("J15").Value = 95
$ws.Range("K15").Value = 90
$ws.Range("W15").Value = 14

# Row 16
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 476
$ws.Range("F16").Value = 258
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 7.79
$ws.Range("I16").Value = 9.210000000000001
$ws.Range("M16").Value = 8

# Row 22
$ws.Range("C22").Value = 22
$ws.Range("D22").Value = 419
$ws.Range("E22").Value = 209
$ws.Range("F22").Value = 210
$ws.Range("G22").Value = 19.05
$ws.Range("H22").Value = 9.5
$ws.Range("I22").Value = 9.550000000000001
$ws.Range("J22").Value = 87
$ws.Range("K22").Value = 95

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 163
$ws.Range("E25").Value = 95
$ws.Range("F25").Value = 68
$ws.Range("G25").Value = 18.11
$ws.Range("H25").Value = 10.56
$ws.Range("I25").Value = 7.56
$ws.Range("J25").Value = 45
$ws.Range("K25").Value = 29
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 2
$ws.Range("V25").Value = 10

# Row 26
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 558
$ws.Range("E26").Value = 238
$ws.Range("F26").Value = 320
$ws.Range("G26").Value = 20.67
$ws.Range("H26").Value = 8.81
$ws.Range("I26").Value = 11.85
$ws.Range("J26").Value = 94
$ws.Range("K26").Value = 100

foreach ($r in 2..26) {
    $ws.Range("AA$r").Value = $timestamp
}
